$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 954, shifting existing rows 954..1036 down to 956..1038
$ws.Rows("954:955").Insert()

# Fill in new row 954
$ws.Cells.Item(954, 1).Value = 5
$ws.Cells.Item(954, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(954, 3).Value = "Maule"
$ws.Cells.Item(954, 4).Value = 45223
$ws.Cells.Item(954, 5).Value = 7
$ws.Cells.Item(954, 6).Value = "Fruta"
$ws.Cells.Item(954, 7).Value = 100102
$ws.Cells.Item(954, 8).Value = "Cítricos"
$ws.Cells.Item(954, 9).Value = 100102005
$ws.Cells.Item(954, 10).Value = "Naranja"
$ws.Cells.Item(954, 11).Value = "Navel Late"
$ws.Cells.Item(954, 12).Value = "Primera"
$ws.Cells.Item(954, 13).Value = 320
$ws.Cells.Item(954, 14).Value = 9000
$ws.Cells.Item(954, 15).Value = 9000
$ws.Cells.Item(954, 16).Value = 9000
$ws.Cells.Item(954, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(954, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(954, 19).Value = 600
$ws.Cells.Item(954, 20).Value = 15

# Fill in new row 955
$ws.Cells.Item(955, 1).Value = 5
$ws.Cells.Item(955, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(955, 3).Value = "Maule"
$ws.Cells.Item(955, 4).Value = 45223
$ws.Cells.Item(955, 5).Value = 7
$ws.Cells.Item(955, 6).Value = "Fruta"
$ws.Cells.Item(955, 7).Value = 100102
$ws.Cells.Item(955, 8).Value = "Cítricos"
$ws.Cells.Item(955, 9).Value = 100102005
$ws.Cells.Item(955, 10).Value = "Naranja"
$ws.Cells.Item(955, 11).Value = "Valencia"
$ws.Cells.Item(955, 12).Value = "Primera"
$ws.Cells.Item(955, 13).Value = 290
$ws.Cells.Item(955, 14).Value = 9000
$ws.Cells.Item(955, 15).Value = 9000
$ws.Cells.Item(955, 16).Value = 9000
$ws.Cells.Item(955, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(955, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(955, 19).Value = 600
$ws.Cells.Item(955, 20).Value = 15
